$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 changes
$ws.Range("D6").Value = 44315
$ws.Range("I6").Value = "Especial"
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 30000
$ws.Range("N6").Value = "$/caja 20 kilos empedrada"
$ws.Range("P6").Value = 1500
$ws.Range("Q6").Value = 20

# Row 7 changes
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("N7").Value = "$/caja 15 kilos granel"
$ws.Range("P7").Value = 1000
$ws.Range("Q7").Value = 15

# Row 8 changes
$ws.Range("D8").Value = 44293
$ws.Range("J8").Value = 10
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 25000
$ws.Range("N8").Value = "$/caja 15 kilos empedrada"
$ws.Range("P8").Value = 1667
